# "new creation changes code"
# - Sheet1!B8 changes from "Fosroc@2" to "Fosroc@1" (selection moves to B8)
# - SEBS_Devloper gets a new row: A3 = "Raj Kumar" (selection moves to A3)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B8").Value = "Fosroc@1"

$ws2 = $wb.Worksheets.Item("SEBS_Devloper")
$ws2.Range("A3").Value = "Raj Kumar"
$ws2.Range("A3").Select()

$ws1.Activate()
$ws1.Range("B8").Select()
